# Rewrite the document's numbering definitions (word/numbering.xml) so
# abstractNumId 1 becomes a standard hybrid-multilevel bullet list
# (numFmt=bullet, Symbol bullet glyph, standard 720/1440/...-twip
# indents) instead of the previous heading-outline "none" list that
# borrowed the Heading1..Heading9 paragraph styles.
#
# We go through Document.Content.WordOpenXML / InsertXML because the
# high-level ListFormat / ListTemplate / ListLevel COM surface only
# knows how to *mint new* list definitions (bumping abstractNumId /
# numId) -- it can't edit an existing <w:abstractNum> in place. Round
# tripping the whole-document Open XML package lets us patch exactly
# the <w:numbering> part while leaving every other part (document.xml,
# styles.xml, ...) untouched.

$d = $word.ActiveDocument

# The new content for word/numbering.xml's single abstractNum (id 1):
# nine levels, all bullets, Symbol font, classic 720/1440/2160/...
# hanging-indent ladder - this is the numbering definition Word mints
# for a default bulleted list.
$newNumbering = @"
<w:numbering xmlns="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:ma="http://schemas.openxmlformats.org/schemaLibrary/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:s="http://schemas.openxmlformats.org/officeDocument/2006/sharedTypes" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:xml="http://www.w3.org/XML/1998/namespace"><w:abstractNum w:abstractNumId="1"><w:multiLevelType w:val="hybridMultilevel"/><w:lvl w:ilvl="0"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:suff w:val="nothing"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Symbol" w:hAnsi="Symbol"/></w:rPr></w:lvl><w:lvl w:ilvl="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:suff w:val="nothing"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Symbol" w:hAnsi="Symbol"/></w:rPr></w:lvl><w:lvl w:ilvl="2"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:suff w:val="nothing"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Symbol" w:hAnsi="Symbol"/></w:rPr></w:lvl><w:lvl w:ilvl="3"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:suff w:val="nothing"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Symbol" w:hAnsi="Symbol"/></w:rPr></w:lvl><w:lvl w:ilvl="4"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:suff w:val="nothing"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Symbol" w:hAnsi="Symbol"/></w:rPr></w:lvl><w:lvl w:ilvl="5"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:suff w:val="nothing"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Symbol" w:hAnsi="Symbol"/></w:rPr></w:lvl><w:lvl w:ilvl="6"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:suff w:val="nothing"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Symbol" w:hAnsi="Symbol"/></w:rPr></w:lvl><w:lvl w:ilvl="7"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:suff w:val="nothing"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Symbol" w:hAnsi="Symbol"/></w:rPr></w:lvl><w:lvl w:ilvl="8"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:suff w:val="nothing"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default" w:ascii="Symbol" w:hAnsi="Symbol"/></w:rPr></w:lvl></w:abstractNum><w:num w:numId="1"><w:abstractNumId w:val="1"/></w:num></w:numbering>
"@

# Pull the whole package (flat OPC) representation of the live
# document so every other part is carried through unchanged.
$package = $d.Content.WordOpenXML

# Swap out just the <pkg:part> that holds word/numbering.xml, keeping
# its pkg:name / pkg:contentType wrapper attributes intact.
$pattern = '(?s)(<pkg:part pkg:name="/word/numbering\.xml"[^>]*><pkg:xmlData>).*?(</pkg:xmlData></pkg:part>)'
$replacement = '${1}' + $newNumbering + '${2}'
$package = [System.Text.RegularExpressions.Regex]::Replace($package, $pattern, $replacement)

# Feed the patched package back in. InsertXML on the whole-document
# range replaces/repopulates every part supplied in the package, so
# the unchanged document.xml / styles.xml / settings.xml content we
# read back from WordOpenXML round-trips as-is while numbering.xml
# picks up the new definition.
$d.Content.InsertXML($package)
